$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("for any SELECT", $false, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null

# Collapse to the end of the found text (end of the "... for any SELECT" paragraph)
$rng.Collapse(0)

# Split off a brand-new paragraph right after the current one, inheriting its
# paragraph formatting (style / numbering / justification).
$rng.InsertParagraphAfter()

# Move into the freshly created paragraph (just past the new paragraph mark)
# and add its text.
$rng.Start = $rng.Start + 1
$rng.End = $rng.Start
$rng.InsertAfter("Save to LLM also the results of each SELECT")
